$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "Java,Python" language value into the new C14 cell
# (row 14 is "Split a Circular Linked List into two halves")
$ws.Range("C14").Value = "Java,Python"

# Update the view: scroll back to top-left A1 and select B7
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
